$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G3").Value = "已售罄"
$ws1.Range("F5").Value = 1191
$ws1.Range("F6").Value = 985
$ws1.Range("F7").Value = 301
$ws1.Range("F10").Value = 921
$ws1.Range("F12").Value = 597
$ws1.Range("F13").Value = 539
$ws1.Range("F14").Value = 1396
$ws1.Range("F15").Value = 128
$ws1.Range("F16").Value = 1305
$ws1.Range("F17").Value = 2962
$ws1.Range("F18").Value = 361
$ws1.Range("F19").Value = 1584
$ws1.Range("F20").Value = 1332
$ws1.Range("F21").Value = 771
$ws1.Range("F22").Value = 222
$ws1.Range("F23").Value = 1321
$ws1.Range("F24").Value = 250
$ws1.Range("F26").Value = 1094
$ws1.Range("F27").Value = 384
$ws1.Range("F28").Value = 3383
$ws1.Range("F29").Value = 658
$ws1.Range("F30").Value = 552
$ws1.Range("F31").Value = 1500

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 17
$ws2.Range("F9").Value = 14
$ws2.Range("F13").Value = 61

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 784
$ws4.Range("F3").Value = 784
$ws4.Range("G4").Value = "已售罄"
$ws4.Range("F7").Value = 1191
$ws4.Range("F8").Value = 985
$ws4.Range("F9").Value = 301
$ws4.Range("F15").Value = 17
$ws4.Range("F19").Value = 14
$ws4.Range("F21").Value = 921
$ws4.Range("F23").Value = 597
$ws4.Range("F24").Value = 539
$ws4.Range("F25").Value = 1396
$ws4.Range("F26").Value = 128
$ws4.Range("F27").Value = 1305
$ws4.Range("F28").Value = 2962
$ws4.Range("F29").Value = 361
$ws4.Range("F30").Value = 1584
$ws4.Range("F31").Value = 1332
$ws4.Range("F32").Value = 771
$ws4.Range("F33").Value = 222
$ws4.Range("F34").Value = 1321
$ws4.Range("F35").Value = 250
$ws4.Range("F39").Value = 1094
$ws4.Range("F40").Value = 384
$ws4.Range("F41").Value = 3383
$ws4.Range("F42").Value = 658
$ws4.Range("F43").Value = 552
$ws4.Range("F44").Value = 1500
$ws4.Range("F45").Value = 61
